$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-02 Saturday" "2024-03-03 Sunday"

Replace-Text "18×91=1638" "87×21=1827"
Replace-Text "78×58=4524" "56×72=4032"
Replace-Text "18×90=1620" "14×22=308"
Replace-Text "42×97=4074" "40×75=3000"
Replace-Text "49×77=3773" "66×15=990"

Replace-Text "62×50=3100" "56×35=1960"
Replace-Text "21×97=2037" "24×67=1608"
Replace-Text "18×79=1422" "59×13=767"
Replace-Text "84×23=1932" "91×34=3094"
Replace-Text "86×12=1032" "72×51=3672"

Replace-Text "85×34=2890" "64×54=3456"
Replace-Text "24×21=504" "34×65=2210"
Replace-Text "50×56=2800" "82×32=2624"
Replace-Text "38×88=3344" "13×79=1027"
Replace-Text "64×56=3584" "31×86=2666"

Replace-Text "76×11=836" "37×14=518"
Replace-Text "47×54=2538" "23×93=2139"
Replace-Text "68×55=3740" "44×75=3300"
Replace-Text "58×81=4698" "33×47=1551"
Replace-Text "37×22=814" "52×50=2600"

Replace-Text "86×29=2494" "20×78=1560"
Replace-Text "13×63=819" "23×87=2001"
Replace-Text "69×25=1725" "77×55=4235"
Replace-Text "77×25=1925" "97×22=2134"
Replace-Text "66×55=3630" "54×67=3618"
